# "Better algorithm when rebordering and clearing needed"
# Update the configuration row-index values on the "Main" sheet (column B)
# that drive the rebordering/clearing algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("B2").Value = 16   # Current Row/Column
$ws.Range("B3").Value = 18   # Wealth Class in Allocation Row
$ws.Range("B4").Value = 6    # Wealth Class in Cash Flow Row
$ws.Range("B5").Value = 7    # Account Row
